$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (values recomputed from new TPM input)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.009821333333333333
$ws.Range("H2").Value = 0.029464
$ws.Range("I2").Value = 0.06297798848338983
$ws.Range("J2").Value = 0.06297798848338984
$ws.Range("M2").Value = 0.003643333333333333
$ws.Range("N2").Value = 0.01093
$ws.Range("O2").Value = 0.002177035403614994
$ws.Range("P2").Value = 0.002177035403614994
$ws.Range("Q2").Value = 0.00003578239111111111
$ws.Range("R2").Value = 0.00032204152
$ws.Range("S2").Value = 0.000137105310576797
$ws.Range("T2").Value = 0.0001371053105767971

# Row 3 updates (values recomputed from new TPM input)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.009821333333333333
$ws.Range("H3").Value = 0.029464
$ws.Range("I3").Value = 0.06297798848338983
$ws.Range("J3").Value = 0.06297798848338984
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.669886333333333
$ws.Range("N3").Value = 5.009659
$ws.Range("O3").Value = 0.997822964596385
$ws.Range("P3").Value = 0.997822964596385
$ws.Range("Q3").Value = 0.01640051030844444
$ws.Range("R3").Value = 0.147604592776
$ws.Range("S3").Value = 0.06284088317281303
$ws.Range("T3").Value = 0.06284088317281304

# Row 4 updates (values recomputed from new TPM input)
$ws.Range("I4").Value = 0.3247949111459754
$ws.Range("J4").Value = 0.3247949111459754
$ws.Range("M4").Value = 0.003643333333333333
$ws.Range("N4").Value = 0.01093
$ws.Range("O4").Value = 0.002177035403614994
$ws.Range("P4").Value = 0.002177035403614994
$ws.Range("Q4").Value = 0.0001845396911111111
$ws.Range("R4").Value = 0.00166085722
$ws.Range("S4").Value = 0.0007070900204787747
$ws.Range("T4").Value = 0.0007070900204787748

# Row 5 updates (values recomputed from new TPM input)
$ws.Range("I5").Value = 0.3247949111459754
$ws.Range("J5").Value = 0.3247949111459754
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.669886333333333
$ws.Range("N5").Value = 5.009659
$ws.Range("O5").Value = 0.997822964596385
$ws.Range("P5").Value = 0.997822964596385
$ws.Range("Q5").Value = 0.08458196929844444
$ws.Range("R5").Value = 0.7612377236860001
$ws.Range("S5").Value = 0.3240878211254966
$ws.Range("T5").Value = 0.3240878211254966

# Row 6 updates (values recomputed from new TPM input)
$ws.Range("G6").Value = 0.09547600000000001
$ws.Range("H6").Value = 0.286428
$ws.Range("I6").Value = 0.6122271003706348
$ws.Range("J6").Value = 0.6122271003706349
$ws.Range("M6").Value = 0.003643333333333333
$ws.Range("N6").Value = 0.01093
$ws.Range("O6").Value = 0.002177035403614994
$ws.Range("P6").Value = 0.002177035403614994
$ws.Range("Q6").Value = 0.0003478508933333334
$ws.Range("R6").Value = 0.00313065804
$ws.Range("S6").Value = 0.001332840072559423
$ws.Range("T6").Value = 0.001332840072559423

# Row 7 updates (values recomputed from new TPM input)
$ws.Range("G7").Value = 0.09547600000000001
$ws.Range("H7").Value = 0.286428
$ws.Range("I7").Value = 0.6122271003706348
$ws.Range("J7").Value = 0.6122271003706349
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.669886333333333
$ws.Range("N7").Value = 5.009659
$ws.Range("O7").Value = 0.997822964596385
$ws.Range("P7").Value = 0.997822964596385
$ws.Range("Q7").Value = 0.1594340675613334
$ws.Range("R7").Value = 1.434906608052
$ws.Range("S7").Value = 0.6108942602980754
$ws.Range("T7").Value = 0.6108942602980755

